# Add more Airbnb test-data rows to the "placesToGo" sheet.
#
# Existing data:
#   A1:C1  location | checkinDate | checkoutDate
#   A2:C2  Seattle  |  12/10/2021 |  12/16/2021
#
# Target data (per the diff):
#   A2:C2  Seattle     |  09/10/2022 |  11/16/2022   (dates updated in place)
#   A3:C3  Los Angeles |  12/12/2021 |  12/25/2021    (new row)
#   A4:C4  Miami       |   9/28/2021 |  11/30/2021    (new row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing Seattle row's check-in / check-out dates.
$ws.Range("B2").Value = " 09/10/2022"
$ws.Range("C2").Value = " 11/16/2022"

# New row 3: Los Angeles. Write cells in this order so new shared-string
# entries are created in the same sequence as the target workbook.
$ws.Range("A3").Value = "Los Angeles"
$ws.Range("B3").Value = " 12/12/2021"

# New row 4: Miami.
$ws.Range("A4").Value = "Miami"
$ws.Range("C3").Value = " 12/25/2021"
$ws.Range("C4").Value = " 11/30/2021"
$ws.Range("B4").Value = " 9/28/2021"

# Copy the date formatting (numFmtId 14) from row 2's date cells onto the
# new rows' date cells instead of re-deriving a number format, so the
# existing style index is reused rather than a new one being minted.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B3:C4").PasteSpecial(-4122) | Out-Null

# Match the workbook's recorded selection after the edits.
$ws.Range("C5").Select() | Out-Null
